$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 44114.285
$ws.Range("J134").Value = 44114.285
$ws.Range("L134").Value = 44114.285
$ws.Range("N134").Value = -54254.285
$ws.Range("H135").Value = 5444.533
$ws.Range("I135").Value = 3848
$ws.Range("J135").Value = 7839.3335
$ws.Range("K135").Value = 34632
$ws.Range("L135").Value = 70554.0015
$ws.Range("M135").Value = -32097
$ws.Range("N135").Value = -75624.0015
$ws.Range("H136").Value = 44950
$ws.Range("J136").Value = 44950
$ws.Range("L136").Value = 44950
$ws.Range("N136").Value = -55150
$ws.Range("H137").Value = 1203.4517
$ws.Range("I137").Value = 880.8
$ws.Range("J137").Value = 1790.091
$ws.Range("K137").Value = 2642.4
$ws.Range("L137").Value = 5370.272999999999
$ws.Range("M137").Value = -92.39999999999964
$ws.Range("N137").Value = -10470.273
$ws.Range("H139").Value = 72420
$ws.Range("J139").Value = 72420
$ws.Range("L139").Value = 72420
$ws.Range("N139").Value = -82700
$ws.Range("H140").Value = 71543.48
$ws.Range("I140").Value = 20000
$ws.Range("J140").Value = 94093.75
$ws.Range("K140").Value = 20000
$ws.Range("L140").Value = 94093.75
$ws.Range("M140").Value = -14820
$ws.Range("N140").Value = -104453.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1275
$ws.Range("I74").Value = 830
$ws.Range("K74").Value = 830
$ws.Range("M74").Value = 44
$ws.Range("H77").Value = 1275
$ws.Range("I77").Value = 830
$ws.Range("K77").Value = 4150
$ws.Range("M77").Value = 218
$ws.Range("H138").Value = 68871.42999999999
$ws.Range("J138").Value = 68871.42999999999
$ws.Range("L138").Value = 68871.42999999999
$ws.Range("N138").Value = -79151.42999999999
$ws.Range("H139").Value = 54722.5
$ws.Range("J139").Value = 54722.5
$ws.Range("L139").Value = 54722.5
$ws.Range("N139").Value = -65002.5
$ws.Range("H140").Value = 105100
$ws.Range("J140").Value = 105100
$ws.Range("L140").Value = 105100
$ws.Range("N140").Value = -115460
$ws.Range("H141").Value = 64600
$ws.Range("J141").Value = 64600
$ws.Range("L141").Value = 64600
$ws.Range("N141").Value = -74960

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 89500
$ws.Range("J140").Value = 89500
$ws.Range("L140").Value = 89500
$ws.Range("N140").Value = -99860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9506.736999999999
$ws.Range("I31").Value = 2972
$ws.Range("J31").Value = 15722.707
$ws.Range("K31").Value = 2972
$ws.Range("L31").Value = 15722.707
$ws.Range("M31").Value = -2677
$ws.Range("N31").Value = -16312.707
$ws.Range("H34").Value = 9506.736999999999
$ws.Range("I34").Value = 2972
$ws.Range("J34").Value = 15722.707
$ws.Range("K34").Value = 2972
$ws.Range("L34").Value = 15722.707
$ws.Range("M34").Value = -2770
$ws.Range("N34").Value = -16126.707
$ws.Range("H58").Value = 1016.5
$ws.Range("I58").Value = 603.65
$ws.Range("J58").Value = 1842.2
$ws.Range("K58").Value = 603.65
$ws.Range("L58").Value = 1842.2
$ws.Range("M58").Value = -400.65
$ws.Range("N58").Value = -2248.2
$ws.Range("H136").Value = 1016.5
$ws.Range("I136").Value = 603.65
$ws.Range("J136").Value = 1842.2
$ws.Range("K136").Value = 1810.95
$ws.Range("L136").Value = 5526.6
$ws.Range("M136").Value = 739.0500000000002
$ws.Range("N136").Value = -10626.6
$ws.Range("H140").Value = 66225
$ws.Range("J140").Value = 66225
$ws.Range("L140").Value = 66225
$ws.Range("N140").Value = -76585

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2900
$ws.Range("I39").Value = 350
$ws.Range("J39").Value = 3218.75
$ws.Range("K39").Value = 1050
$ws.Range("L39").Value = 9656.25
$ws.Range("M39").Value = -756
$ws.Range("N39").Value = -10244.25
$ws.Range("H131").Value = 5155522
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 6411110.5
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 19233331.5
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -19243411.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H132").Value = 2276.7144
$ws.Range("I132").Value = 1928.4688
$ws.Range("J132").Value = 3391.1
$ws.Range("K132").Value = 5785.4064
$ws.Range("L132").Value = 10173.3
$ws.Range("M132").Value = -3255.4064
$ws.Range("N132").Value = -15233.3
$ws.Range("H138").Value = 69071.42999999999
$ws.Range("J138").Value = 69071.42999999999
$ws.Range("L138").Value = 69071.42999999999
$ws.Range("N138").Value = -79351.42999999999
$ws.Range("H140").Value = 89864.5
$ws.Range("J140").Value = 89864.5
$ws.Range("L140").Value = 89864.5
$ws.Range("N140").Value = -100224.5
$ws.Range("H141").Value = 53000
$ws.Range("J141").Value = 53000
$ws.Range("L141").Value = 53000
$ws.Range("N141").Value = -63360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 451.5
$ws.Range("I19").Value = 451.5
$ws.Range("K19").Value = 451.5
$ws.Range("M19").Value = -281.5
$ws.Range("H132").Value = 4750.2085
$ws.Range("I132").Value = 4856.95
$ws.Range("K132").Value = 14570.85
$ws.Range("M132").Value = -12040.85
$ws.Range("H136").Value = 4026.8125
$ws.Range("I136").Value = 4500.8887
$ws.Range("J136").Value = 3417.2856
$ws.Range("K136").Value = 13502.6661
$ws.Range("L136").Value = 10251.8568
$ws.Range("M136").Value = -10952.6661
$ws.Range("N136").Value = -15351.8568
$ws.Range("H138").Value = 58623.363
$ws.Range("J138").Value = 58623.363
$ws.Range("L138").Value = 58623.363
$ws.Range("N138").Value = -68903.363
$ws.Range("H139").Value = 57216.668
$ws.Range("J139").Value = 67660
$ws.Range("L139").Value = 67660
$ws.Range("N139").Value = -77940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 1550
$ws.Range("I20").Value = 1550
$ws.Range("K20").Value = 1550
$ws.Range("M20").Value = -1310
$ws.Range("H138").Value = 69366.664
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 69366.664
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 69366.664
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -79646.664
$ws.Range("H139").Value = 58116.668
$ws.Range("J139").Value = 58116.668
$ws.Range("L139").Value = 58116.668
$ws.Range("N139").Value = -68396.66800000001
$ws.Range("H140").Value = 29950
$ws.Range("J140").Value = 29950
$ws.Range("L140").Value = 29950
$ws.Range("N140").Value = -40310
$ws.Range("H141").Value = 79014.375
$ws.Range("J141").Value = 79014.375
$ws.Range("L141").Value = 79014.375
$ws.Range("N141").Value = -89374.375
